$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'322.80"
$ws.Range("E2").Value = "'-2.64%"
$ws.Range("D3").Value = "'43.07"
$ws.Range("E3").Value = "'-5.20%"
$ws.Range("D4").Value = "'5.182"
$ws.Range("E4").Value = "'-7.60%"
$ws.Range("D5").Value = "'0.08203"
$ws.Range("E5").Value = "'-1.84%"
$ws.Range("D6").Value = "'4.322"
$ws.Range("E6").Value = "'-2.66%"
$ws.Range("D7").Value = "'1.834"
$ws.Range("E7").Value = "'-11.18%"
$ws.Range("D8").Value = "'0.9339"
$ws.Range("E8").Value = "'-2.99%"
$ws.Range("E9").Value = "'-4.86%"
$ws.Range("D10").Value = "'0.1871"
$ws.Range("E10").Value = "'-2.64%"
$ws.Range("D11").Value = "'0.09438"
$ws.Range("E11").Value = "'-4.52%"
$ws.Range("D12").Value = "'0.04626"
$ws.Range("E12").Value = "'0.19%"
$ws.Range("D13").Value = "'7.429"
$ws.Range("E13").Value = "'-28.82%"
$ws.Range("E14").Value = "'-0.34%"
$ws.Range("D15").Value = "'0.001298"
$ws.Range("E15").Value = "'1.03%"
$ws.Range("D16").Value = "'0.005790"
$ws.Range("E16").Value = "'-4.81%"
$ws.Range("D17").Value = "'3.362"
$ws.Range("E17").Value = "'-0.45%"
$ws.Range("D18").Value = "'2.540"
$ws.Range("E18").Value = "'-1.66%"
$ws.Range("D19").Value = "'0.3376"
$ws.Range("E19").Value = "'0.20%"
$ws.Range("D20").Value = "'0.1389"
$ws.Range("E20").Value = "'-0.31%"
$ws.Range("E21").Value = "'-1.27%"
$ws.Range("D22").Value = "'0.04161"
$ws.Range("E22").Value = "'-0.62%"
$ws.Range("D23").Value = "'0.001246"
$ws.Range("E23").Value = "'-5.29%"
$ws.Range("D24").Value = "'0.004339"
$ws.Range("E24").Value = "'-4.98%"
$ws.Range("E25").Value = "'-8.06%"
$ws.Range("D26").Value = "'0.0002979"
$ws.Range("E26").Value = "'-20.57%"
$ws.Range("D38").Value = "'0.02769"
$ws.Range("E38").Value = "'2.03%"
$ws.Range("D39").Value = "'0.05582"
$ws.Range("E39").Value = "'-3.19%"
$ws.Range("D40").Value = "'0.008139"
$ws.Range("E40").Value = "'3.71%"
$ws.Range("E41").Value = "'-2.55%"
$ws.Range("D42").Value = "'0.006544"
$ws.Range("E42").Value = "'-10.33%"
$ws.Range("D43").Value = "'0.002091"
$ws.Range("E43").Value = "'3.69%"
$ws.Range("D44").Value = "'0.007515"
$ws.Range("E44").Value = "'-17.33%"
$ws.Range("D45").Value = "'0.3487"
$ws.Range("E45").Value = "'-1.76%"
$ws.Range("D46").Value = "'0.00006982"
$ws.Range("E46").Value = "'-2.39%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.27%"
$ws.Range("D48").Value = "'0.003480"
$ws.Range("E48").Value = "'-0.53%"
$ws.Range("D49").Value = "'0.003530"
$ws.Range("E49").Value = "'0.64%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.27%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.27%"
